# Reading planner update for April
# Fills in the "Bible Chapter" (column B) values for rows 92-121,
# which correspond to the dates 2023-04-01 through 2023-04-30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "1 KINGS 1-4",
    "1 KINGS 1-4",
    "1 KINGS 5-8",
    "1 KINGS 9-11",
    "1 KINGS 12-16",
    "1 KINGS 17-19",
    "1 KINGS 20-22",
    "2 KINGS 1-3",
    "2 KINGS 1-3",
    "2 KINGS 4-8",
    "2 KINGS 9-12",
    "2 KINGS 13-17",
    "2 KINGS 18-21",
    "2 KINGS 22-25",
    "1 CHRONICLES 1-9",
    "1 CHRONICLES 1-9",
    "1 CHRONICLES 10-16",
    "1 CHRONICLES 17-21",
    "1 CHRONICLES 22-27",
    "1 CHRONICLES 28-29",
    "2 CHRONICLES 1-5",
    "2 CHRONICLES 6-9",
    "2 CHRONICLES 6-9",
    "2 CHRONICLES 10-12",
    "2 CHRONICLES 13-16",
    "2 CHRONICLES 17-20",
    "2 CHRONICLES 21-25",
    "2 CHRONICLES 26-28",
    "2 CHRONICLES 29-32",
    "2 CHRONICLES 33-36"
)

$startRow = 92
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Update the visible window / active selection to reflect where the
# user was working (matches the sheetView change in the diff).
$ws.Range("B122").Select()
$excel.ActiveWindow.ScrollRow = 116
$excel.ActiveWindow.ScrollColumn = 1
